# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.331.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.589.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '190.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.634'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.586.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.663'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("E13").Value = '  +2.95%  '
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.164.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.583.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.164.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '478.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '92.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("E27").Value = '  -1.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("E32").Value = '  +3.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '584.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.31%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0801'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.398'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +22.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.140'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.46%  '
$ws.Range("E42").Value = '  -5.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.241.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("E44").Value = '  +7.16%  '
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0445'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.10%  '
